$wb = $excel.ActiveWorkbook
$students = $wb.Worksheets.Item("Students")

# Add "Errors" sheet right after "Students"
$errorsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $students)
$errorsSheet.Name = "Errors"

# Add "Warnings" sheet right after "Errors"
$warningsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $errorsSheet)
$warningsSheet.Name = "Warnings"

# Error messages reported for the "Students" sheet, row 2
$messages = @(
    "''Sheet ""Students"" Row: 2 Missing ""LAST NAME""',",
    "''Sheet ""Students"" Row: 2 Missing ""FIRST NAME""',",
    "''Sheet ""Students"" Row: 2 Missing ""STUDENT ID""',",
    "''Sheet ""Students"" Row: 2 Missing ""BIRTH DT""',",
    "''Sheet ""Students"" Row: 2 Missing ""OFF CLS""',",
    "''Sheet ""Students"" Row: 2 Invalid birthday """"',"
)

for ($i = 0; $i -lt $messages.Length; $i++) {
    $row = $i + 1
    $cell = $errorsSheet.Range("A$row")
    $cell.Value = $messages[$i]
    # Drop the implicit "quote prefix" text style Excel applies whenever a
    # value is typed starting with an apostrophe - the message itself
    # should carry the literal leading quote as plain text.
    $cell.Style = "Normal"
}

$null = $errorsSheet.Range("A1:A6").Select()

# Restore the original selection on "Students" and make sure it is no
# longer the tab shown when the workbook opens.
$null = $students.Select()
$null = $students.Range("G2").Select()

# "Warnings" is the sheet that should be active/selected when the
# workbook is opened.
$null = $warningsSheet.Select()
